$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L (Nb nouveaux décès à l'hôpital) and M (Nb nouveaux décès
# extra-hospitaliers) are formatted as Text ("@"), so assigning a numeric
# Value directly would store it as a text string. Briefly switch the cell
# to a General number format, write the number, then restore the original
# Text format so the stored value stays numeric (matching the target file).
function Set-NumericOnTextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "General"
    $rng.Value = $value
    $rng.NumberFormat = "@"
}

# Row 617: update number of new positive cases
$ws.Range("C617").Value = 53

# Row 618: update new cases, ICU patients, and new in-hospital deaths
$ws.Range("C618").Value = 72
$ws.Range("G618").Value = 11
Set-NumericOnTextCell "L618" 1

# Row 619: fill in previously empty daily figures
$ws.Range("C619").Value = 58
$ws.Range("E619").Value = 3
$ws.Range("F619").Value = 3
$ws.Range("G619").Value = 8
Set-NumericOnTextCell "L619" 0
Set-NumericOnTextCell "M619" 0

# Row 620: fill in previously empty daily figures
$ws.Range("C620").Value = 69
$ws.Range("E620").Value = 3
$ws.Range("F620").Value = 3
$ws.Range("G620").Value = 8
Set-NumericOnTextCell "L620" 0
Set-NumericOnTextCell "M620" 0

# Row 621: fill in previously empty daily figures
$ws.Range("C621").Value = 29
$ws.Range("E621").Value = 3
$ws.Range("F621").Value = 3
$ws.Range("G621").Value = 10
Set-NumericOnTextCell "L621" 0
Set-NumericOnTextCell "M621" 0

# Row 622: fill in previously empty daily figures
$ws.Range("C622").Value = 4
$ws.Range("E622").Value = 4
$ws.Range("F622").Value = 3
$ws.Range("G622").Value = 13
Set-NumericOnTextCell "L622" 0
Set-NumericOnTextCell "M622" 0
